# Update country COVID-19 statistics and re-sort by total cases (descending).
# Source data for this update (27 March 2020, 17:14): country name plus
# Casos totales / Nuevos casos / Casos activos / Recuperados / Casos criticos /
# Muertes hoy / Muertes, indexed by destination row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$countryRows = @(
    @(4,"Estados Unidos",93131,7696,2417,89332,2432,87,1382),
    @(5,"China",81340,55,74588,3460,1034,5,3292),
    @(6,"Italia",80589,0,10361,62013,3612,0,8215),
    @(7,"España",64059,6273,9357,49768,4165,569,4934),
    @(8,"Alemania",49344,5406,5673,43367,23,37,304),
    @(9,"Iran",32332,2926,11133,18821,2893,144,2378),
    @(10,"Francia",29155,0,4948,22511,3375,0,1696),
    @(11,"Reino Unido",14543,2885,135,13649,163,181,759),
    @(12,"Suiza",12311,500,897,11207,203,15,207),
    @(13,"Corea del Sur",9332,91,4528,4665,59,8,139),
    @(14,"Paises Bajos",8603,1172,3,8054,761,112,546),
    @(15,"Austria",7500,591,225,7217,128,9,58),
    @(16,"Belgica",7284,1049,858,6137,690,69,289),
    @(17,"Portugal",4268,724,43,4149,71,16,76),
    @(18,"Canada",4043,0,228,3776,120,0,39),
    @(19,"Noruega",3694,322,6,3671,70,3,17),
    @(20,"Turquia",3629,0,26,3528,136,0,75),
    @(21,"Australia",3180,130,170,2997,23,0,13),
    @(22,"Suecia",3046,206,16,2938,214,15,92),
    @(23,"Israel",3035,342,79,2944,49,4,12),
    @(24,"Brasil",3027,42,6,2944,296,0,77),
    @(25,"Malasia",2161,130,259,1876,54,3,26),
    @(26,"Chequia",2062,137,11,2042,34,0,9),
    @(27,"Dinamarca",2010,133,1,1957,109,11,52),
    @(28,"Irlanda",1819,0,5,1795,47,0,19),
    @(29,"Chile",1610,304,43,1562,7,1,5),
    @(30,"Luxemburgo",1605,152,40,1550,25,6,15),
    @(31,"Ecuador",1403,0,3,1366,58,0,34),
    @(32,"Japon",1387,0,359,981,57,0,47),
    @(33,"Pakistan",1296,95,23,1264,7,0,9),
    @(34,"Rumania",1292,263,115,1153,32,1,24),
    @(35,"Polonia",1289,68,7,1266,3,0,16),
    @(36,"Tailandia",1136,91,97,1034,11,1,5),
    @(37,"Arabia Saudita",1104,92,35,1066,6,0,3),
    @(38,"Indonesia",1046,153,46,913,0,9,87),
    @(39,"Finlandia",1041,83,10,1024,32,2,7),
    @(40,"Rusia",1036,196,45,988,8,0,3),
    @(41,"Sudafrica",927,0,12,913,7,2,2),
    @(42,"Grecia",892,0,42,823,53,0,27),
    @(43,"Islandia",890,88,97,791,18,0,2),
    @(44,"India",874,147,73,781,0,0,20),
    @(45,"Filipinas",803,96,31,718,1,9,54),
    @(46,"Singapur",732,49,183,547,17,0,2),
    @(47,"Crucero",712,0,597,105,15,0,10),
    @(48,"Panama",674,0,2,663,20,0,9),
    @(49,"Eslovenia",632,70,10,613,14,3,9),
    @(50,"Argentina",589,0,72,504,0,1,13),
    @(51,"Croacia",586,91,37,546,14,0,3),
    @(52,"Mexico",585,110,4,573,1,2,8),
    @(53,"Republica Dominicana",581,93,3,558,0,10,20),
    @(54,"Peru",580,0,14,557,14,0,9),
    @(55,"Estonia",575,37,11,563,7,0,1),
    @(56,"Catar",549,0,43,506,6,0,0),
    @(57,"Serbia",528,71,15,505,25,1,8),
    @(58,"Hong Kong",518,64,111,403,5,0,4),
    @(59,"Egipto",495,0,102,369,0,0,24),
    @(60,"Colombia",491,0,8,477,0,0,6),
    @(61,"Barein",466,8,227,235,1,0,4),
    @(62,"Irak",458,76,122,296,0,4,40),
    @(63,"Libano",391,23,23,361,3,1,7),
    @(64,"Nueva Zelanda",368,0,37,331,1,0,0),
    @(65,"Argelia",367,0,29,313,0,0,25),
    @(66,"Lituania",345,46,1,339,2,1,5),
    @(67,"Emiratos Arabes Unidos",333,0,52,279,2,0,2),
    @(68,"Armenia",329,39,28,300,6,0,1),
    @(69,"Hungria",300,39,34,256,6,0,10),
    @(70,"Bulgaria",293,29,9,281,8,0,3),
    @(71,"Letonia",280,36,1,279,0,0,0),
    @(72,"Marruecos",275,0,8,256,1,0,11),
    @(73,"Eslovaquia",269,43,2,267,1,0,0),
    @(74,"Principado de Andorra",267,43,1,263,11,0,3),
    @(75,"Taiwan",267,15,30,235,0,0,2),
    @(76,"Uruguay",238,0,0,238,3,0,0),
    @(77,"Costa Rica",231,0,2,227,5,0,2),
    @(78,"Bosnia y Herzegovina",231,40,5,222,1,1,4),
    @(79,"Tunez",227,30,2,219,10,1,6),
    @(80,"Ucrania",226,30,5,216,0,0,5),
    @(81,"Kuwait",225,17,57,168,11,0,0),
    @(82,"San Marino",223,15,4,198,12,0,21),
    @(83,"Republica de Macedonia",219,18,3,213,1,0,3),
    @(84,"Jordania",212,0,2,210,0,0,0),
    @(85,"Moldavia",199,22,2,195,33,1,2),
    @(86,"Albania",186,12,31,147,3,2,8),
    @(87,"Burkina Faso",180,28,12,159,0,2,9),
    @(88,"Azerbaiyan",165,43,15,147,6,0,3),
    @(89,"Vietnam",163,10,20,143,3,0,0),
    @(90,"Republica de Chipre",146,0,4,139,3,0,3),
    @(91,"Reunion",145,10,1,144,0,0,0),
    @(92,"Islas Feroe",144,4,54,90,2,0,0),
    @(93,"Malta",139,5,2,137,1,0,0),
    @(94,"Kazajistan",137,24,3,133,0,0,1),
    @(95,"Ghana",136,4,1,131,1,0,4),
    @(96,"Oman",131,22,23,108,0,0,0),
    @(97,"Senegal",119,14,11,108,0,0,0),
    @(98,"Brunei",115,1,11,104,1,0,0),
    @(99,"Venezuela",107,0,31,75,2,0,1),
    @(100,"Sri Lanka",106,0,7,99,5,0,0),
    @(101,"Camboya",99,1,11,88,1,0,0),
    @(102,"Costa de Marfil",96,0,3,93,0,0,0),
    @(103,"Mauricio",94,13,0,92,1,0,2),
    @(104,"Afganistan",94,0,2,88,0,0,4),
    @(105,"Bielorrusia",94,8,32,62,2,0,0),
    @(106,"Estado de Palestina",91,5,17,73,0,0,1),
    @(107,"Camerun",88,13,2,84,0,1,2),
    @(108,"Uzbekistan",88,13,5,82,8,1,1),
    @(109,"Martinica",81,0,0,80,12,0,1),
    @(110,"Georgia",81,2,13,68,1,0,0),
    @(111,"Cuba",80,13,4,74,2,0,2),
    @(112,"Guadalupe",73,0,0,72,4,0,1),
    @(113,"Montenegro",70,1,0,69,1,0,1),
    @(114,"Honduras",68,1,0,67,0,0,1),
    @(115,"Trinidad yTobago",66,1,1,63,0,1,2),
    @(116,"Nigeria",65,0,3,61,0,0,1),
    @(117,"Bolivia",61,0,0,61,0,0,0),
    @(118,"Kirguistan",58,14,0,58,0,0,0),
    @(119,"Liechtenstein",56,0,0,56,0,0,0),
    @(120,"Gibraltar",55,20,14,41,0,0,0),
    @(121,"Paraguay",52,11,1,48,1,0,3),
    @(122,"Consejo Danes para los Refugiados",51,0,2,46,0,0,3),
    @(123,"Mayotte",50,14,0,50,0,0,0),
    @(124,"Ruanda",50,0,0,50,0,0,0),
    @(125,"Banglades",48,4,11,32,1,0,5),
    @(126,"Puerto Rico",39,0,1,36,0,0,2),
    @(127,"Macao",34,1,10,24,0,0,0),
    @(128,"Monaco",33,0,1,32,0,0,0),
    @(129,"Guam",32,0,0,31,0,0,1),
    @(130,"Kenia",31,0,1,29,0,0,1),
    @(131,"Polinesia Francesa",30,0,0,30,0,0,0),
    @(132,"Isla de Man",29,3,0,29,0,0,0),
    @(133,"Aruba",28,0,1,27,0,0,0),
    @(134,"Guayana Francesa",28,0,6,22,0,0,0),
    @(135,"Jamaica",26,0,2,23,0,0,1),
    @(136,"Togo",25,1,1,24,0,0,0),
    @(137,"Guatemala",25,0,4,20,0,0,1),
    @(138,"Barbados",24,0,0,24,0,0,0),
    @(139,"Madagascar",24,1,0,24,0,0,0),
    @(140,"Zambia",22,6,0,22,0,0,0),
    @(141,"Uganda",18,4,0,18,0,0,0),
    @(142,"Islas Virgenes de los Estados Unidos",17,0,0,17,0,0,0),
    @(143,"Etiopia",16,4,0,16,0,0,0),
    @(144,"Nueva Caledonia",15,1,0,15,0,0,0),
    @(145,"Bermudas",15,0,2,13,0,0,0),
    @(146,"Maldivas",14,1,9,5,0,0,0),
    @(147,"El Salvador",13,0,0,13,0,0,0),
    @(148,"Tanzania",13,0,1,12,0,0,0),
    @(149,"Republica de Yibuti",12,1,0,12,0,0,0),
    @(150,"Guinea Ecuatorial",12,0,0,12,0,0,0),
    @(151,"Mongolia",11,0,0,11,0,0,0),
    @(152,"Mali",11,7,0,11,0,0,0),
    @(153,"Dominica",11,0,0,11,0,0,0),
    @(154,"San Martin (Parte Francesa)",11,0,0,11,0,0,0),
    @(155,"Niger",10,0,0,9,0,0,1),
    @(156,"Groenlandia",10,4,2,8,0,0,0),
    @(157,"Bahamas",9,0,1,8,0,0,0),
    @(158,"Surinam",8,0,0,8,0,0,0),
    @(159,"Guinea",8,4,0,8,0,0,0),
    @(160,"Haiti",8,0,0,8,0,0,0),
    @(161,"Islas Caimanes",8,0,0,7,0,0,1),
    @(162,"Namibia",8,0,2,6,0,0,0),
    @(163,"Antigua y Barbuda",7,0,0,7,0,0,0),
    @(164,"Seychelles",7,0,0,7,0,0,0),
    @(165,"Mozambique",7,0,0,7,0,0,0),
    @(166,"Granada",7,0,0,7,0,0,0),
    @(167,"Gabon",7,0,0,6,0,0,1),
    @(168,"Curazao",7,0,2,4,0,0,1),
    @(169,"Laos",6,0,0,6,0,0,0),
    @(170,"Benin",6,0,0,6,0,0,0),
    @(171,"Eritrea",6,0,0,6,0,0,0),
    @(172,"Suazilandia",6,0,0,6,0,0,0),
    @(173,"Montserrat",5,0,0,5,0,0,0),
    @(174,"Siria",5,0,0,5,0,0,0),
    @(175,"Birmania",5,0,0,5,0,0,0),
    @(176,"Fiyi",5,0,0,5,0,0,0),
    @(177,"Guyana",5,0,0,4,0,0,1),
    @(178,"Cabo Verde",5,0,0,4,0,0,1),
    @(179,"Zimbabue",5,2,0,4,0,0,1),
    @(180,"Angola",4,0,0,4,0,0,0),
    @(181,"Santa Sede",4,0,0,4,0,0,0),
    @(182,"Congo",4,0,0,4,0,0,0),
    @(183,"Nepal",4,1,1,3,0,0,0),
    @(184,"Republica de Africa Central",3,0,0,3,0,0,0),
    @(185,"Republica del Chad",3,0,0,3,0,0,0),
    @(186,"Butan",3,1,0,3,0,0,0),
    @(187,"Liberia",3,0,0,3,0,0,0),
    @(188,"San Bartolome",3,0,0,3,0,0,0),
    @(189,"San Martin (Parte Holandesa)",3,0,0,3,0,0,0),
    @(190,"Somalia",3,1,0,3,0,0,0),
    @(191,"Mauritania",3,0,0,3,0,0,0),
    @(192,"Sudan",3,0,0,2,0,0,1),
    @(193,"Santa Lucia",3,0,1,2,0,0,0),
    @(194,"Gambia",3,0,0,2,0,0,1),
    @(195,"Anguila",2,0,0,2,0,0,0),
    @(196,"Islas Turcas y Caicos",2,0,0,2,0,0,0),
    @(197,"Belice",2,0,0,2,0,0,0),
    @(198,"Guinea-Bisau",2,0,0,2,0,0,0),
    @(199,"San Cristobal y Nieves",2,0,0,2,0,0,0),
    @(200,"Islas Virgenes Britanicas",2,0,0,2,0,0,0),
    @(201,"Nicaragua",2,0,0,1,0,1,1),
    @(202,"Timor Oriental",1,0,0,1,0,0,0),
    @(203,"San Vicente y las Granadinas",1,0,0,1,0,0,0),
    @(204,"Papua Nueva Guinea",1,0,0,1,0,0,0),
    @(205,"Libia",1,0,0,1,0,0,0)
)

foreach ($row in $countryRows) {
    $targetRow = $row[0]
    $ws.Cells.Item($targetRow, 1).Value = $row[1]
    $ws.Cells.Item($targetRow, 2).Value = $row[2]
    $ws.Cells.Item($targetRow, 3).Value = $row[3]
    $ws.Cells.Item($targetRow, 4).Value = $row[4]
    $ws.Cells.Item($targetRow, 5).Value = $row[5]
    $ws.Cells.Item($targetRow, 6).Value = $row[6]
    $ws.Cells.Item($targetRow, 7).Value = $row[7]
    $ws.Cells.Item($targetRow, 8).Value = $row[8]
}

# Update the "last refreshed" timestamp shown above the table.
$ws.Range("A1").Value = "Datos actualizados a 27 de Marzo de 2020 a las 17:14"

Write-Host "Updated $($countryRows.Count) country rows and refreshed timestamp."
